# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 290 (pushing the existing rows
# 290-350 down to 291-351) for "Femacal de La Calera" / Acelga.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 290:350 down to 291:351, leaving row 290 empty for the new entry.
$ws.Rows("290:290").Insert()

$ws.Range("A290").Value = 3
$ws.Range("B290").Value = "Femacal de La Calera"
$ws.Range("C290").Value = "Coquimbo"
$ws.Range("D290").Value2 = 44711
$ws.Range("E290").Value = 5
$ws.Range("F290").Value = 100112009
$ws.Range("G290").Value = "Acelga"
$ws.Range("H290").Value = "Sin especificar"
$ws.Range("I290").Value = "Primera"
$ws.Range("J290").Value = 270
$ws.Range("K290").Value = 3000
$ws.Range("L290").Value = 3300
$ws.Range("M290").Value = 3133
$ws.Range("N290").Value = "`$/docena de atados (6 kilos)"
$ws.Range("O290").Value = "Provincia de Quillota"
$ws.Range("P290").Value = 522
$ws.Range("Q290").Value = 6
$ws.Range("R290").Value = "Hortaliza"
